$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

$rng = $ws.Range("A1:F38")
$rng.NumberFormat = "@"

$ws.Range("A1").Value = "iuliia.1"
$ws.Range("B1").Value = "CAD"
$ws.Range("C1").Value = "287.52"
$ws.Range("D1").Value = "14"
$ws.Range("E1").Value = "Visa"
$ws.Range("F1").Value = "3371001216"

$ws.Range("A2").Value = "iuliia.4"
$ws.Range("B2").Value = "AED"
$ws.Range("C2").Value = "735.25"
$ws.Range("D2").Value = "40"
$ws.Range("E2").Value = "MasterCard"
$ws.Range("F2").Value = "3347847168"

$ws.Range("A3").Value = "iuliia.4"
$ws.Range("B3").Value = "KWD"
$ws.Range("C3").Value = "63.42"
$ws.Range("D3").Value = "4"
$ws.Range("E3").Value = "American Express"
$ws.Range("F3").Value = "3308438016"

$ws.Range("A4").Value = "iuliia.2"
$ws.Range("B4").Value = "NZD"
$ws.Range("C4").Value = "2105.35"
$ws.Range("D4").Value = "200"
$ws.Range("E4").Value = "Visa"
$ws.Range("F4").Value = "3325735680"

$ws.Range("A5").Value = "iuliia.6"
$ws.Range("B5").Value = "AUD"
$ws.Range("C5").Value = "5790.55"
$ws.Range("D5").Value = "160"
$ws.Range("E5").Value = "MasterCard"
$ws.Range("F5").Value = "3365369600"

$ws.Range("A6").Value = "iuliia.6"
$ws.Range("B6").Value = "GBP"
$ws.Range("C6").Value = "2654.76"
$ws.Range("D6").Value = "125"
$ws.Range("E6").Value = "American Express"
$ws.Range("F6").Value = "3308594432"

$ws.Range("A7").Value = "iuliia.1"
$ws.Range("B7").Value = "CAD"
$ws.Range("C7").Value = "301.52"
$ws.Range("D7").Value = "14"
$ws.Range("E7").Value = "Visa"
$ws.Range("F7").ClearContents()

$ws.Range("A8").Value = "iuliia.4"
$ws.Range("B8").Value = "AED"
$ws.Range("C8").Value = "824.08"
$ws.Range("D8").Value = "40"
$ws.Range("E8").Value = "MasterCard"
$ws.Range("F8").ClearContents()

$ws.Range("A9").Value = "iuliia.4"
$ws.Range("B9").Value = "KWD"
$ws.Range("C9").Value = "67.42"
$ws.Range("D9").Value = "4"
$ws.Range("E9").Value = "American Express"
$ws.Range("F9").ClearContents()

$ws.Range("A10").Value = "iuliia.2"
$ws.Range("B10").Value = "NZD"
$ws.Range("C10").Value = "2313.64"
$ws.Range("D10").Value = "200"
$ws.Range("E10").Value = "Visa"
$ws.Range("F10").ClearContents()

$ws.Range("A11").Value = "iuliia.6"
$ws.Range("B11").Value = "AUD"
$ws.Range("C11").Value = "6200.83"
$ws.Range("D11").Value = "160"
$ws.Range("E11").Value = "MasterCard"
$ws.Range("F11").ClearContents()

$ws.Range("A12").Value = "iuliia.6"
$ws.Range("B12").Value = "GBP"
$ws.Range("C12").Value = "2767.41"
$ws.Range("D12").Value = "125"
$ws.Range("E12").Value = "American Express"
$ws.Range("F12").ClearContents()

$ws.Range("A13").Value = "iuliia.1"
$ws.Range("B13").Value = "CAD"
$ws.Range("C13").Value = "315.52"
$ws.Range("D13").Value = "14"
$ws.Range("E13").Value = "Visa"
$ws.Range("F13").ClearContents()

$ws.Range("A14").Value = "iuliia.1"
$ws.Range("B14").Value = "CAD"
$ws.Range("C14").Value = "329.52"
$ws.Range("D14").Value = "14"
$ws.Range("E14").Value = "Visa"
$ws.Range("F14").Value = "3348859392"

$ws.Range("A15").Value = "iuliia.4"
$ws.Range("B15").Value = "AED"
$ws.Range("C15").Value = "824.08"
$ws.Range("D15").Value = "40"
$ws.Range("E15").Value = "MasterCard"
$ws.Range("F15").ClearContents()

$ws.Range("A16").Value = "iuliia.1"
$ws.Range("B16").Value = "CAD"
$ws.Range("C16").Value = "329.52"
$ws.Range("D16").Value = "14"
$ws.Range("E16").Value = "Visa"
$ws.Range("F16").ClearContents()

$ws.Range("A17").Value = "iuliia.1"
$ws.Range("B17").Value = "CAD"
$ws.Range("C17").Value = "357.52"
$ws.Range("D17").Value = "14"
$ws.Range("E17").Value = "Visa"
$ws.Range("F17").ClearContents()

$ws.Range("A18").Value = "iuliia.1"
$ws.Range("B18").Value = "CAD"
$ws.Range("C18").Value = "357.52"
$ws.Range("D18").Value = "14"
$ws.Range("E18").Value = "Visa"
$ws.Range("F18").ClearContents()

$ws.Range("A19").Value = "iuliia.4"
$ws.Range("B19").Value = "AED"
$ws.Range("C19").Value = "824.08"
$ws.Range("D19").Value = "40"
$ws.Range("E19").Value = "MasterCard"
$ws.Range("F19").ClearContents()

$ws.Range("A20").Value = "iuliia.1"
$ws.Range("B20").Value = "CAD"
$ws.Range("C20").Value = "371.52"
$ws.Range("D20").Value = "14"
$ws.Range("E20").Value = "Visa"
$ws.Range("F20").Value = "3376005339"

$ws.Range("A21").Value = "iuliia.4"
$ws.Range("B21").Value = "AED"
$ws.Range("C21").Value = "824.08"
$ws.Range("D21").Value = "40"
$ws.Range("E21").Value = "MasterCard"
$ws.Range("F21").Value = "3306281175"

$ws.Range("A22").Value = "iuliia.4"
$ws.Range("B22").Value = "KWD"
$ws.Range("C22").Value = "70.69"
$ws.Range("D22").Value = "4"
$ws.Range("E22").Value = "American Express"
$ws.Range("F22").Value = "3373696546"

$ws.Range("A23").Value = "iuliia.2"
$ws.Range("B23").Value = "NZD"
$ws.Range("C23").Value = "2313.64"
$ws.Range("D23").Value = "200"
$ws.Range("E23").Value = "Visa"
$ws.Range("F23").Value = "3318716612"

$ws.Range("A24").Value = "iuliia.4"
$ws.Range("B24").Value = "KWD"
$ws.Range("C24").Value = "74.68"
$ws.Range("D24").Value = "4"
$ws.Range("E24").Value = "American Express"
$ws.Range("F24").Value = "3388898189"

$ws.Range("A25").Value = "iuliia.2"
$ws.Range("B25").Value = "NZD"
$ws.Range("C25").Value = "2521.94"
$ws.Range("D25").Value = "200"
$ws.Range("E25").Value = "Visa"
$ws.Range("F25").Value = "3353686220"

$ws.Range("A26").Value = "iuliia.6"
$ws.Range("B26").Value = "AUD"
$ws.Range("C26").Value = "6200.83"
$ws.Range("D26").Value = "160"
$ws.Range("E26").Value = "MasterCard"
$ws.Range("F26").Value = "3355506230"

$ws.Range("A27").Value = "iuliia.6"
$ws.Range("B27").Value = "AUD"
$ws.Range("C27").Value = "6358.69"
$ws.Range("D27").Value = "160"
$ws.Range("E27").Value = "MasterCard"
$ws.Range("F27").Value = "3340162017"

$ws.Range("A28").Value = "iuliia.6"
$ws.Range("B28").Value = "GBP"
$ws.Range("C28").Value = "2908.32"
$ws.Range("D28").Value = "125"
$ws.Range("E28").Value = "American Express"
$ws.Range("F28").Value = "3357642259"

$ws.Range("A29").Value = "iuliia.1"
$ws.Range("B29").Value = "CAD"
$ws.Range("C29").Value = "399.52"
$ws.Range("D29").Value = "14"
$ws.Range("E29").Value = "Visa"
$ws.Range("F29").ClearContents()

$ws.Range("A30").Value = "iuliia.1"
$ws.Range("B30").Value = "CAD"
$ws.Range("C30").Value = "413.52"
$ws.Range("D30").Value = "14"
$ws.Range("E30").Value = "Visa"
$ws.Range("F30").ClearContents()

$ws.Range("A31").Value = "iuliia.4"
$ws.Range("B31").Value = "AED"
$ws.Range("C31").Value = "961.74"
$ws.Range("D31").Value = "40"
$ws.Range("E31").Value = "MasterCard"
$ws.Range("F31").ClearContents()

$ws.Range("A32").Value = "iuliia.1"
$ws.Range("B32").Value = "CAD"
$ws.Range("C32").Value = "427.52"
$ws.Range("D32").Value = "14"
$ws.Range("E32").Value = "Visa"
$ws.Range("F32").ClearContents()

$ws.Range("A33").Value = "iuliia.1"
$ws.Range("B33").Value = "CAD"
$ws.Range("C33").Value = "427.52"
$ws.Range("D33").Value = "14"
$ws.Range("E33").Value = "Visa"
$ws.Range("F33").ClearContents()

$ws.Range("A34").Value = "iuliia.4"
$ws.Range("B34").Value = "AED"
$ws.Range("C34").Value = "961.74"
$ws.Range("D34").Value = "40"
$ws.Range("E34").Value = "MasterCard"
$ws.Range("F34").Value = "3341297037"

$ws.Range("A35").Value = "iuliia.4"
$ws.Range("B35").Value = "KWD"
$ws.Range("C35").Value = "81.95"
$ws.Range("D35").Value = "4"
$ws.Range("E35").Value = "American Express"
$ws.Range("F35").Value = "3337332298"

$ws.Range("A36").Value = "iuliia.2"
$ws.Range("B36").Value = "NZD"
$ws.Range("C36").Value = "2730.25"
$ws.Range("D36").Value = "200"
$ws.Range("E36").Value = "Visa"
$ws.Range("F36").Value = "3336089182"

$ws.Range("A37").Value = "iuliia.6"
$ws.Range("B37").Value = "AUD"
$ws.Range("C37").Value = "6768.98"
$ws.Range("D37").Value = "160"
$ws.Range("E37").Value = "MasterCard"
$ws.Range("F37").Value = "3309999836"

$ws.Range("A38").Value = "iuliia.6"
$ws.Range("B38").Value = "GBP"
$ws.Range("C38").Value = "3020.98"
$ws.Range("D38").Value = "125"
$ws.Range("E38").Value = "American Express"
$ws.Range("F38").Value = "3370312014"

$rng.Style = "Normal"